$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-12 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-13 Saturday", 2) | Out-Null
$d.Content.Find.Execute("17×65=1105", $true, $false, $false, $false, $false, $true, 1, $false, "79×67=5293", 2) | Out-Null
$d.Content.Find.Execute("28×85=2380", $true, $false, $false, $false, $false, $true, 1, $false, "52×66=3432", 2) | Out-Null
$d.Content.Find.Execute("84×55=4620", $true, $false, $false, $false, $false, $true, 1, $false, "70×44=3080", 2) | Out-Null
$d.Content.Find.Execute("39×74=2886", $true, $false, $false, $false, $false, $true, 1, $false, "34×50=1700", 2) | Out-Null
$d.Content.Find.Execute("47×56=2632", $true, $false, $false, $false, $false, $true, 1, $false, "70×36=2520", 2) | Out-Null
$d.Content.Find.Execute("37×62=2294", $true, $false, $false, $false, $false, $true, 1, $false, "14×25=350", 2) | Out-Null
$d.Content.Find.Execute("15×92=1380", $true, $false, $false, $false, $false, $true, 1, $false, "33×79=2607", 2) | Out-Null
$d.Content.Find.Execute("22×92=2024", $true, $false, $false, $false, $false, $true, 1, $false, "20×47=940", 2) | Out-Null
$d.Content.Find.Execute("16×58=928", $true, $false, $false, $false, $false, $true, 1, $false, "32×57=1824", 2) | Out-Null
$d.Content.Find.Execute("31×84=2604", $true, $false, $false, $false, $false, $true, 1, $false, "12×90=1080", 2) | Out-Null
$d.Content.Find.Execute("72×21=1512", $true, $false, $false, $false, $false, $true, 1, $false, "19×24=456", 2) | Out-Null
$d.Content.Find.Execute("16×45=720", $true, $false, $false, $false, $false, $true, 1, $false, "76×24=1824", 2) | Out-Null
$d.Content.Find.Execute("29×97=2813", $true, $false, $false, $false, $false, $true, 1, $false, "36×80=2880", 2) | Out-Null
$d.Content.Find.Execute("41×37=1517", $true, $false, $false, $false, $false, $true, 1, $false, "97×88=8536", 2) | Out-Null
$d.Content.Find.Execute("52×18=936", $true, $false, $false, $false, $false, $true, 1, $false, "65×25=1625", 2) | Out-Null
$d.Content.Find.Execute("66×58=3828", $true, $false, $false, $false, $false, $true, 1, $false, "35×91=3185", 2) | Out-Null
$d.Content.Find.Execute("39×47=1833", $true, $false, $false, $false, $false, $true, 1, $false, "72×53=3816", 2) | Out-Null
$d.Content.Find.Execute("32×95=3040", $true, $false, $false, $false, $false, $true, 1, $false, "81×40=3240", 2) | Out-Null
$d.Content.Find.Execute("20×69=1380", $true, $false, $false, $false, $false, $true, 1, $false, "94×96=9024", 2) | Out-Null
$d.Content.Find.Execute("33×18=594", $true, $false, $false, $false, $false, $true, 1, $false, "63×21=1323", 2) | Out-Null
$d.Content.Find.Execute("54×56=3024", $true, $false, $false, $false, $false, $true, 1, $false, "45×33=1485", 2) | Out-Null
$d.Content.Find.Execute("25×30=750", $true, $false, $false, $false, $false, $true, 1, $false, "23×77=1771", 2) | Out-Null
$d.Content.Find.Execute("34×47=1598", $true, $false, $false, $false, $false, $true, 1, $false, "71×90=6390", 2) | Out-Null
$d.Content.Find.Execute("27×26=702", $true, $false, $false, $false, $false, $true, 1, $false, "77×22=1694", 2) | Out-Null
$d.Content.Find.Execute("88×92=8096", $true, $false, $false, $false, $false, $true, 1, $false, "53×99=5247", 2) | Out-Null
